# Apply the UI-design-driven changes to the Holidays template workbook:
#   1. Rename the "type" column header (E1) to "holiday_type".
#   2. Move the active cell / selection from G10 to E7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header in E1 from "type" to "holiday_type"
$ws.Range("E1").Value = "holiday_type"

# 2. Update the selected / active cell to E7
[void]$ws.Range("E7").Select()
